$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("validateCreateCustomerApi")
$ws.Activate()

# New customer ids pulled in after adding the java mail api.
# The first three overwrite the values already sitting in A13:A15,
# the rest extend the "id" list further down through A36.
$ids = @(
    "cus_Lj1sjwM4YV9p2I",
    "cus_Lj1sOAiI4mUdRr",
    "cus_Lj1sYclhwbI5jQ",
    "cus_Lj1ryzYh8gGfCX",
    "cus_Lj1rmfOczqlhF5",
    "cus_Lj1r5jc4SiZbY2",
    "cus_LixxC5zDMvX4hp",
    "cus_LixxSxX3HA3WQP",
    "cus_LixxfBxSPabjHI",
    "cus_LigOyXqqT5j9ls",
    "cus_LigObgiOZVWTzZ",
    "cus_LigOPUIj43ZELc",
    "cus_LifMUzfnGKmBoz",
    "cus_LifM6nLEZCLKq7",
    "cus_LifMlrqjwRKPRZ",
    "cus_LifESzNhWEB74Y",
    "cus_LifENDKWOcz3qH",
    "cus_LifERMKVrWRC9q",
    "cus_Lif4IrNFXq3yyn",
    "cus_Lif4jHeB4HZUpJ",
    "cus_Lif4hz73M9gdtF",
    "cus_LibCeM6Nmdmq40",
    "cus_LibC7Yc31XkWZp",
    "cus_LibC6ACVRYrJrP"
)

$startRow = 13
for ($i = 0; $i -lt $ids.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $ids[$i]
}

# Match the author's saved selection (A13:A36, anchored at A13) and scroll
# the window down so row 16 is the first visible row.
$ws.Range("A13:A36").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
